# Update "想去人数" (want-to-go count) values in the F column on both the
# "展览" sheet and the aggregated "全部类型" sheet, matching the regenerated
# site data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 531
    "F8"  = 49
    "F9"  = 3971
    "F10" = 4301
    "F12" = 132
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
